# Merge #225 - #241
# Remove the "Menu" (F) and "Program" (G) columns from Sheet1's header/data
# rows 1-2 (and their now-orphaned shared-string values "aa"/"swqd"), then
# switch the active sheet/selection from the "Note" sheet back to "Sheet1".

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Delete columns F:G (Menu, Program) entirely - shifts H.. left to F..
$ws1.Range("F1:G1").EntireColumn.Delete()

# Note sheet keeps its own cell selection, just no longer the active tab.
$ws2.Range("L6").Select()

# Sheet1 becomes the active/selected sheet & tab, with K8 selected.
$ws1.Activate()
$ws1.Range("K8").Select()
